# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row in the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 12:47:20"
$wsZhCn.Range("H2").Value = "2016-03-18 12:47:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 12:47:23"
$wsDeDe.Range("H2").Value = "2016-03-18 12:47:42"
